$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column C (UserHourCost, formerly "Cost") formatting into new column D
# (entire data range incl. header) so that D inherits the same per-row styles.
$ws.Range("C1:C38").Copy()
$ws.Range("D1:D38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rename / re-order the header labels:
#   A = UserCostPodr   (was UserPodr)
#   B = CostUserName   (was UserName)
#   C = UserHourCost   (was Cost)
#   D = UserMonthCost  (new)
$ws.Range("A1").Value = "UserCostPodr"
$ws.Range("B1").Value = "CostUserName"
$ws.Range("C1").Value = "UserHourCost"
$ws.Range("D1").Value = "UserMonthCost"

# Column widths: widen column A (drop its old best-fit) and extend the
# best-fit-style width used by column C to the new column D as well.
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(4).ColumnWidth = 19

# Update the stored sort state so it spans B:D instead of B:C.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add2($ws.Range("B2:B24"))
$ws.Sort.SetRange($ws.Range("B2:D24"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()
